# Generate Report for Handoff
#
# Updates the localization-status report:
#  - Status text "Handed back: in sync with en-US" -> "Ready for handoff"
#    (Overview!E2:F2, zh-cn!C2, de-de!C2)
#  - "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
#    bumped to reflect the newly generated handoff package
#  - Narrower Status-column widths on all three sheets (status text got shorter)

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$ws_overview.Range("E2").Value = "Ready for handoff"
$ws_overview.Range("F2").Value = "Ready for handoff"
$ws_zhcn.Range("C2").Value = "Ready for handoff"
$ws_dede.Range("C2").Value = "Ready for handoff"

# --- Timestamps ---
$ws_overview.Range("G2").Value = "2016-09-06 05:16:09"
$ws_dede.Range("H2").Value = "2016-09-06 05:16:09"
$ws_zhcn.Range("H2").Value = "2016-09-06 05:15:58"

# --- Column width: shrink the Status column(s) from ~30 chars to ~17 chars ---
# (the host quantizes ColumnWidth to 1/6-character steps, so 16.333333333333332
# is the input that reproducibly yields the narrowest achievable stored width)
$ws_overview.Range("E:E").ColumnWidth = 16.333333333333332
$ws_overview.Range("F:F").ColumnWidth = 16.333333333333332
$ws_zhcn.Range("C:C").ColumnWidth = 16.333333333333332
$ws_dede.Range("C:C").ColumnWidth = 16.333333333333332
